$d = $word.ActiveDocument

# Replace "Prog:" with "File:" everywhere it occurs (both text-box copies)
$d.Content.Find.Execute("Prog:", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "File:", 2)

Write-Output "done"
